$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.988.89"
$ws.Range("E2").Value = "'  -2.08%  "
$ws.Range("D3").Value = "'1.985.31"
$ws.Range("E3").Value = "'  -1.32%  "
$ws.Range("D4").Value = "'1.016"
$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("D5").Value = "'329.76"
$ws.Range("E5").Value = "'  -0.75%  "
$ws.Range("D6").Value = "'1.014"
$ws.Range("E6").Value = "'  +0.15%  "
$ws.Range("D7").Value = "'0.4921"
$ws.Range("E7").Value = "'  -2.36%  "
$ws.Range("D8").Value = "'0.4160"
$ws.Range("E8").Value = "'  -2.29%  "
$ws.Range("D9").Value = "'55.03"
$ws.Range("E9").Value = "'  +2.05%  "
$ws.Range("D10").Value = "'0.08814"
$ws.Range("E10").Value = "'  -4.54%  "
$ws.Range("D11").Value = "'1.082"
$ws.Range("E11").Value = "'  -4.00%  "
$ws.Range("D12").Value = "'2.069.08"
$ws.Range("E12").Value = "'  +5.45%  "
$ws.Range("D13").Value = "'22.70"
$ws.Range("E13").Value = "'  -3.79%  "
$ws.Range("D14").Value = "'7.907"
$ws.Range("E14").Value = "'  -2.75%  "
$ws.Range("D15").Value = "'6.366"
$ws.Range("E15").Value = "'  -2.95%  "
$ws.Range("D16").Value = "'1.017"
$ws.Range("E16").Value = "'  +0.30%  "
$ws.Range("D17").Value = "'91.52"
$ws.Range("E17").Value = "'  -4.74%  "
$ws.Range("D18").Value = "'0.00001097"
$ws.Range("E18").Value = "'  -2.50%  "
$ws.Range("D19").Value = "'0.06664"
$ws.Range("E19").Value = "'  -0.21%  "
$ws.Range("D20").Value = "'19.26"
$ws.Range("E20").Value = "'  -3.84%  "
$ws.Range("D21").Value = "'1.013"
$ws.Range("E21").Value = "'  +0.16%  "
$ws.Range("D22").Value = "'5.944"
$ws.Range("E22").Value = "'  -0.95%  "
$ws.Range("D23").Value = "'29.041.68"
$ws.Range("E23").Value = "'  -2.08%  "
$ws.Range("D24").Value = "'11.87"
$ws.Range("E24").Value = "'  -1.20%  "
$ws.Range("D25").Value = "'2.313"
$ws.Range("E25").Value = "'  +1.12%  "
$ws.Range("D26").Value = "'2.290.40"
$ws.Range("E26").Value = "'  +3.85%  "
$ws.Range("D27").Value = "'20.70"
$ws.Range("D28").Value = "'156.85"
$ws.Range("E28").Value = "'  -1.82%  "
$ws.Range("D29").Value = "'6.191"
$ws.Range("E29").Value = "'  -4.06%  "
$ws.Range("D30").Value = "'2.219"
$ws.Range("E30").Value = "'  -5.69%  "
$ws.Range("D31").Value = "'126.08"
$ws.Range("E31").Value = "'  -2.00%  "
$ws.Range("D32").Value = "'1.034"
$ws.Range("E32").Value = "'  -2.64%  "
$ws.Range("D33").Value = "'0.09827"
$ws.Range("E33").Value = "'  -1.49%  "
$ws.Range("D34").Value = "'1.512"
$ws.Range("E34").Value = "'  -5.32%  "
$ws.Range("D35").Value = "'5.810"
$ws.Range("E35").Value = "'  -1.39%  "
$ws.Range("D36").Value = "'3.738"
$ws.Range("E36").Value = "'  -1.67%  "
$ws.Range("D37").Value = "'0.02397"
$ws.Range("E37").Value = "'  -3.21%  "
$ws.Range("D38").Value = "'1.302"
$ws.Range("E38").Value = "'  -2.04%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06355"
$ws.Range("E39").Value = "'  -0.87%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.968"
$ws.Range("E40").Value = "'  -6.74%  "
$ws.Range("D41").Value = "'0.6436"
$ws.Range("E41").Value = "'  -2.41%  "
$ws.Range("D42").Value = "'11.44"
$ws.Range("E42").Value = "'  -3.42%  "
$ws.Range("E43").Value = "'  -5.50%  "
$ws.Range("D44").Value = "'1.014"
$ws.Range("E44").Value = "'  +0.20%  "
$ws.Range("D45").Value = "'1.358"
$ws.Range("E45").Value = "'  +5.45%  "
$ws.Range("D46").Value = "'0.6143"
$ws.Range("E46").Value = "'  -3.72%  "
$ws.Range("D47").Value = "'13.28"
$ws.Range("E47").Value = "'  -3.09%  "
$ws.Range("D48").Value = "'2.147"
$ws.Range("E48").Value = "'  -3.21%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000346"
$ws.Range("E49").Value = "'  +7.37%  "
$ws.Range("B50").Value = "PancakeSwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D50").Value = "'3.480"
$ws.Range("E50").Value = "'  -1.71%  "
$ws.Range("D51").Value = "'2.156"
$ws.Range("E51").Value = "'  +5.82%  "
